# Applies the cryptos list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.828.48"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.640.03"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.82"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.30%  "
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").Value = "1.876.30"
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "1.640.90"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("E14").Value = "  +3.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.45"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.45%  "
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "29.852.50"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.96"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.72%  "
$ws.Range("D20").Value = "0.0₃0703"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.93"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.54"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.109"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.63"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0494"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.20"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("D34").Value = "1.422.19"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.69"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.72"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.35%  "
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.46"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.562"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.833"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.783.36"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.36"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "49.36"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -9.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "93.19"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.80%  "
$ws.Range("E51").Value = "  +0.66%  "

Write-Host "Applied 86 cell updates to cryptos sheet."
